$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "{%p if jobs.there_are_any %}" -> "{%p if jobs.count > 0 %}"
#    and move the "_GoBack" bookmark from the child.there_are_any paragraph
#    into this line (right after "count > 0").
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("jobs.there_are_any", $true, $false, $false, $false, $false, $true, 1, $false, "jobs.count > 0", 2)

# locate the freshly inserted "count > 0" text so we know where to drop the bookmark
$countRange = $d.Content
$countRange.Find.Execute("count > 0", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# remove the bookmark from its old location (if present) and re-add it here
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

$bookmarkSpot = $d.Range($countRange.End, $countRange.End)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot)

# ---------------------------------------------------------------------------
# 2) "{%p if income" + "s" -> "{%p if incomes" (merge the two runs)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("incomes.there_are_any", $true, $false, $false, $false, $false, $true, 1, $false, "incomes.there_are_any", 2)

# ---------------------------------------------------------------------------
# 3) "{%p for " + "inc" + " in " + "incomes" + " %}" -> "{%p for inc in incomes %}"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("{%p for inc in incomes %}", $true, $false, $false, $false, $false, $true, 1, $false, "{%p for inc in incomes %}", 2)

# ---------------------------------------------------------------------------
# 4) "{{" + "inc" + "." + "description" + "}}" -> "{{inc.description}}"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("{{inc.description}}", $true, $false, $false, $false, $false, $true, 1, $false, "{{inc.description}}", 2)
